$wb = $excel.ActiveWorkbook

# --- Rename sheet "Historiadelapoblación de Asi" -> "HistoriadelapoblacióndeAsi" ---
$wsHist = $wb.Worksheets.Item("Historiadelapoblación de Asi")
$wsHist.Name = "HistoriadelapoblacióndeAsi"

# --- Update the chart series formulas on GraficaHistoria so they reference the
#     renamed sheet without the (now unnecessary) quoting ---
$wsGraf = $wb.Worksheets.Item("GraficaHistoria")
$co = $wsGraf.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection()
$s1 = $series.Item(1)
$s1.Formula = "=SERIES(HistoriadelapoblacióndeAsi!`$A`$76,HistoriadelapoblacióndeAsi!`$A`$2:`$A`$76,HistoriadelapoblacióndeAsi!`$A`$2:`$A`$76,1)"
$s2 = $series.Item(2)
$s2.Formula = "=SERIES(HistoriadelapoblacióndeAsi!`$B`$1,,HistoriadelapoblacióndeAsi!`$B`$2:`$B`$76,2)"

# --- PaisesdeAsia: scroll the view back to the top-left (drops topLeftCell="A42"),
#     keep the existing selection ---
$wsPaises = $wb.Worksheets.Item("PaisesdeAsia")
$wsPaises.Activate()
$wsPaises.Range("C65").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# --- HistoriadelapoblacióndeAsi becomes the active/selected tab with a new
#     selected cell F22 ---
$wsHist.Activate()
$wsHist.Range("F22").Select()
